# Introduction_metagenomics.pptx — "Added intro to genome assembly"
#
# Slide 20, shape "object 2" (the big bullet-list text box) lists the
# reference-based-assembly software options. The three one-word bullets
# "Newbler" / "AMOS" / "MIRA" are replaced by a single bullet listing
# "Spade [ metaspede, metavir etc]", so the paragraph count for that
# sub-list drops from 3 to 1 and the auto-fit text box shrinks.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(20)
$shape = $s.Shapes.Item(2)
$tf = $shape.TextFrame
$tr = $tf.TextRange

# Paragraph 3 currently reads "Newbler" — turn it into the new bullet,
# built up run-by-run (mirrors how the author's edit ended up with one
# run per typed/autocorrected word) while preserving the paragraph's
# existing run formatting (sz 2200, spc -5, Arial MT).
$newbler = $tr.Paragraphs(3, 1)
$newbler.Text = "Spade [ "
$newbler.InsertAfter("metaspede") | Out-Null
$newbler.InsertAfter(", ") | Out-Null
$newbler.InsertAfter("metavir") | Out-Null
$newbler.InsertAfter(" ") | Out-Null
$newbler.InsertAfter("etc") | Out-Null
$newbler.InsertAfter("]") | Out-Null

# Paragraphs 4 ("AMOS") and 5 ("MIRA") are now redundant — delete both,
# paragraph marks included, collapsing the sub-list down to one bullet.
$tr2 = $tf.TextRange
$amos = $tr2.Paragraphs(4, 1)
$mira = $tr2.Paragraphs(5, 1)
$start = $amos.Start
$len = ($mira.Start + $mira.Length) - $start
$tr2.Characters($start, $len).Delete()

# The text box has spAutoFit; removing two paragraphs of lvl-2 text
# reflows it shorter. Apply the resulting laid-out height explicitly.
$shape.Height = 5094985 / 12700.0
